$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1117.2285
$ws.Range("I40").Value = 1063.1904
$ws.Range("J40").Value = 1198.2858
$ws.Range("K40").Value = 1063.1904
$ws.Range("L40").Value = 1198.2858
$ws.Range("M40").Value = -888.1904
$ws.Range("N40").Value = -1548.2858
$ws.Range("H43").Value = 33468.445
$ws.Range("I43").Value = 47999.25
$ws.Range("J43").Value = 21843.8
$ws.Range("K43").Value = 47999.25
$ws.Range("L43").Value = 21843.8
$ws.Range("M43").Value = -47930.25
$ws.Range("N43").Value = -21981.8
$ws.Range("H55").Value = 468.57144
$ws.Range("I55").Value = 317.1
$ws.Range("J55").Value = 847.25
$ws.Range("K55").Value = 317.1
$ws.Range("L55").Value = 847.25
$ws.Range("M55").Value = -103.1
$ws.Range("N55").Value = -1275.25
$ws.Range("H62").Value = 74893
$ws.Range("I62").Value = 132286.88
$ws.Range("J62").Value = 9300
$ws.Range("K62").Value = 132286.88
$ws.Range("L62").Value = 9300
$ws.Range("M62").Value = -131662.88
$ws.Range("N62").Value = -10548
$ws.Range("H65").Value = 74893
$ws.Range("I65").Value = 132286.88
$ws.Range("J65").Value = 9300
$ws.Range("K65").Value = 661434.4
$ws.Range("L65").Value = 46500
$ws.Range("M65").Value = -658314.4
$ws.Range("N65").Value = -52740
$ws.Range("J100").Value = 3400
$ws.Range("L100").Value = 3400
$ws.Range("N100").Value = -4482
$ws.Range("H134").Value = 114999.5
$ws.Range("J134").Value = 114999.5
$ws.Range("L134").Value = 114999.5
$ws.Range("N134").Value = -125139.5
$ws.Range("H135").Value = 2097.818
$ws.Range("I135").Value = 1907.6
$ws.Range("K135").Value = 17168.4
$ws.Range("M135").Value = -14633.4
$ws.Range("H138").Value = 3586.89
$ws.Range("I138").Value = 1445.2354
$ws.Range("J138").Value = 4025.5422
$ws.Range("K138").Value = 4335.706200000001
$ws.Range("L138").Value = 12076.6266
$ws.Range("M138").Value = 804.2937999999995
$ws.Range("N138").Value = -22356.6266

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2081.2727
$ws.Range("I122").Value = 1943.7778
$ws.Range("K122").Value = 5831.3334
$ws.Range("M122").Value = -3381.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1414.85
$ws.Range("I107").Value = 1372
$ws.Range("J107").Value = 1657.6666
$ws.Range("K107").Value = 1372
$ws.Range("L107").Value = 1657.6666
$ws.Range("M107").Value = 548
$ws.Range("N107").Value = -5497.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3333.3333
$ws.Range("I16").Value = 3333.3333
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3333.3333
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3046.3333
$ws.Range("N16").Value = ""
$ws.Range("H31").Value = 4531.28
$ws.Range("I31").Value = 2505.5293
$ws.Range("J31").Value = 8836
$ws.Range("K31").Value = 2505.5293
$ws.Range("L31").Value = 8836
$ws.Range("M31").Value = -2210.5293
$ws.Range("N31").Value = -9426
$ws.Range("H34").Value = 4531.28
$ws.Range("I34").Value = 2505.5293
$ws.Range("J34").Value = 8836
$ws.Range("K34").Value = 2505.5293
$ws.Range("L34").Value = 8836
$ws.Range("M34").Value = -2303.5293
$ws.Range("N34").Value = -9240
$ws.Range("H58").Value = 2910.7778
$ws.Range("I58").Value = 3279.6
$ws.Range("K58").Value = 3279.6
$ws.Range("M58").Value = -3076.6
$ws.Range("H99").Value = 4876.6665
$ws.Range("J99").Value = 5500
$ws.Range("L99").Value = 5500
$ws.Range("N99").Value = -8496
$ws.Range("H107").Value = 766.03705
$ws.Range("I107").Value = 829.2632
$ws.Range("K107").Value = 829.2632
$ws.Range("M107").Value = 1090.7368
$ws.Range("H113").Value = 3333.3333
$ws.Range("I113").Value = 3333.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3333.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1163.3333
$ws.Range("N113").Value = ""
$ws.Range("H126").Value = 4876.6665
$ws.Range("J126").Value = 5500
$ws.Range("L126").Value = 16500
$ws.Range("N126").Value = -21440
$ws.Range("H134").Value = 15160.678
$ws.Range("I134").Value = 8175.222
$ws.Range("J134").Value = 62312.5
$ws.Range("K134").Value = 24525.666
$ws.Range("L134").Value = 186937.5
$ws.Range("M134").Value = -21990.666
$ws.Range("N134").Value = -192007.5
$ws.Range("H136").Value = 2910.7778
$ws.Range("I136").Value = 3279.6
$ws.Range("K136").Value = 9838.8
$ws.Range("M136").Value = -7288.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 87.73913
$ws.Range("I23").Value = 96.666664
$ws.Range("K23").Value = 289.999992
$ws.Range("M23").Value = -54.99999200000002
$ws.Range("H40").Value = 223.31818
$ws.Range("I40").Value = 149.71428
$ws.Range("J40").Value = 352.125
$ws.Range("K40").Value = 598.85712
$ws.Range("L40").Value = 1408.5
$ws.Range("M40").Value = -529.85712
$ws.Range("N40").Value = -1546.5
$ws.Range("H86").Value = 415.875
$ws.Range("I86").Value = 360.75
$ws.Range("K86").Value = 1082.25
$ws.Range("M86").Value = 103.75
$ws.Range("H89").Value = 415.875
$ws.Range("I89").Value = 360.75
$ws.Range("K89").Value = 3246.75
$ws.Range("M89").Value = 2681.25
$ws.Range("H129").Value = 535.5833
$ws.Range("I129").Value = 535.5833
$ws.Range("K129").Value = 1606.7499
$ws.Range("M129").Value = 3393.2501
$ws.Range("H131").Value = 63970.86
$ws.Range("I131").Value = 400365.8
$ws.Range("J131").Value = 9713.613
$ws.Range("K131").Value = 1201097.4
$ws.Range("L131").Value = 29140.839
$ws.Range("M131").Value = -1196057.4
$ws.Range("N131").Value = -39220.839
$ws.Range("H137").Value = 4030.2942
$ws.Range("I137").Value = 2484.6667
$ws.Range("K137").Value = 7454.000100000001
$ws.Range("M137").Value = -2354.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19167.334
$ws.Range("I7").Value = 21200.8
$ws.Range("J7").Value = 9000
$ws.Range("K7").Value = 21200.8
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = -21088.8
$ws.Range("N7").Value = -9224
$ws.Range("H16").Value = 13161774
$ws.Range("I16").Value = 22728540
$ws.Range("K16").Value = 22728540
$ws.Range("M16").Value = -22728370
$ws.Range("H22").Value = 925.6667
$ws.Range("J22").Value = 1097
$ws.Range("L22").Value = 1097
$ws.Range("N22").Value = -1687
$ws.Range("H27").Value = 925.6667
$ws.Range("J27").Value = 1097
$ws.Range("L27").Value = 1097
$ws.Range("N27").Value = -1311
$ws.Range("H40").Value = 4426.846
$ws.Range("I40").Value = 4034.6
$ws.Range("K40").Value = 4034.6
$ws.Range("M40").Value = -3898.6
$ws.Range("H46").Value = 1814.5
$ws.Range("J46").Value = 1666.7
$ws.Range("L46").Value = 1666.7
$ws.Range("N46").Value = -2042.7
$ws.Range("H55").Value = 199.88889
$ws.Range("J55").Value = 298.5
$ws.Range("L55").Value = 298.5
$ws.Range("N55").Value = -644.5
$ws.Range("H68").Value = 12530
$ws.Range("I68").Value = 4225
$ws.Range("J68").Value = 18066.666
$ws.Range("K68").Value = 4225
$ws.Range("L68").Value = 18066.666
$ws.Range("M68").Value = -3476
$ws.Range("N68").Value = -19564.666
$ws.Range("H71").Value = 12530
$ws.Range("I71").Value = 4225
$ws.Range("J71").Value = 18066.666
$ws.Range("K71").Value = 21125
$ws.Range("L71").Value = 90333.33
$ws.Range("M71").Value = -17381
$ws.Range("N71").Value = -97821.33
$ws.Range("H82").Value = 2426.2778
$ws.Range("J82").Value = 2625.3333
$ws.Range("L82").Value = 2625.3333
$ws.Range("N82").Value = -3347.3333
$ws.Range("H85").Value = 2426.2778
$ws.Range("J85").Value = 2625.3333
$ws.Range("L85").Value = 2625.3333
$ws.Range("N85").Value = -5121.3333
$ws.Range("H93").Value = 446179.8
$ws.Range("I93").Value = 557446.75
$ws.Range("J93").Value = 1112
$ws.Range("K93").Value = 557446.75
$ws.Range("L93").Value = 1112
$ws.Range("M93").Value = -556198.75
$ws.Range("N93").Value = -3608
$ws.Range("H126").Value = 19167.334
$ws.Range("I126").Value = 21200.8
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 63602.39999999999
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -61132.39999999999
$ws.Range("N126").Value = -31940
$ws.Range("H132").Value = 5628.143
$ws.Range("I132").Value = 4611.875
$ws.Range("J132").Value = 6983.1665
$ws.Range("K132").Value = 13835.625
$ws.Range("L132").Value = 20949.4995
$ws.Range("M132").Value = -11305.625
$ws.Range("N132").Value = -26009.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 510
$ws.Range("I113").Value = 183.33333
$ws.Range("K113").Value = 549.99999
$ws.Range("M113").Value = 1620.00001
$ws.Range("H126").Value = 4977.5
$ws.Range("I126").Value = 4721.875
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 14165.625
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -11695.625
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 5939.5366
$ws.Range("I132").Value = 3338.853
$ws.Range("K132").Value = 10016.559
$ws.Range("M132").Value = -7486.559000000001

